$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" -----------
# Overview sheet tracks status per-locale in columns E (zh-cn) and F (de-de).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-locale detail sheets carry the same status in column C ("Status").
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Shrink the columns that held the longer status string -----------------
# The shorter "In Translation" text no longer needs as much room, so the
# status columns are re-sized to fit the new content.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
